$d = $word.ActiveDocument

# Locate the paragraph that ends with the "Gitignore" heading text
# (Chapter-2 ... Gitignore) so the new explanatory paragraph can be
# inserted directly after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Gitignore\s*\r?$") {
        $target = $p
    }
}

if ($target -eq $null) {
    $target = $d.Paragraphs.Last
}

# Create a brand new (empty) paragraph right after the heading.
$target.Range.InsertParagraphAfter()

# Re-fetch it via its position (Next of the heading paragraph) so we
# operate on the freshly inserted paragraph.
$newPara = $target.Next()

$paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>Gitignore is a file which is made by the programmer to ignore certain files/certain types of files/certain directories from being tracked by git. It is made in the git repository which is being tracked by git. It is a text file with the name of &#8220;.gitignore&#8221;. Inside it, the name of all files which need to be ignored, the type of files which need to ignored (*.&lt;file extension&gt;), directories which need to be ignored (/&lt;name of dir&gt;) for the single directory which needs to be ignored, or if there are multiple directories</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> with same name and all of them need to be ignored then (&lt;name of dir need to be ignored/&gt;).</w:t></w:r></w:p>'

$newPara.Range.InsertXML($paraXml)

Write-Output "Inserted Gitignore description paragraph."
